# Update the weekly Fruta/Hortaliza price records: the rows of data (D, L,
# M, N, O, P, Q, R, S, T columns) are shuffled to new row positions, as if
# re-sorted/re-ordered by date. Capture the current per-row values first
# (since several target rows are also source rows), then write them back to
# their new destinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together for each data row (2..13).
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot current values for every data row before writing anything.
# (.Value2 is used for reading - .Value's getter surfaces a COM variant
# wrapper rather than the scalar in this runtime.)
$snapshot = @{}
for ($r = 2; $r -le 13; $r++) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $row
}

# Destination row -> source row (which old row's data now lives here).
$mapping = @{
    2  = 7
    3  = 13
    4  = 5
    5  = 2
    6  = 9
    7  = 10
    8  = 6
    9  = 12
    10 = 3
    11 = 8
    12 = 4
    13 = 11
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcData[$c]
    }
}
